$wb = $excel.ActiveWorkbook

# Date/time number format shared by column A across all sheets.
$dateFmt = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------
# Sheet 1 "ROW50-FE-LIFTER" - append new row 93
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(93, 1).Value = 45770.80108277778
$ws1.Cells.Item(93, 1).NumberFormat = $dateFmt
$ws1.Cells.Item(93, 2).Value = "0x01,0x90"
$ws1.Cells.Item(93, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item(93, 4).Value = "0x01,0x3e"
$ws1.Cells.Item(93, 5).Value = "0xe"
$ws1.Cells.Item(93, 6).Value = 400
$ws1.Cells.Item(93, 7).Value = [double]"5.68631262647114e+23"
$ws1.Cells.Item(93, 8).Value = 318
$ws1.Cells.Item(93, 9).Value = 14

# ---------------------------------------------------------------
# Sheet 2 "ROW50-MID-LIFTER" - append new row 95
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(95, 1).Value = 45770.76311342593
$ws2.Cells.Item(95, 1).NumberFormat = $dateFmt
$ws2.Cells.Item(95, 2).Value = "0x01,0x90 "
$ws2.Cells.Item(95, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item(95, 4).Value = "0x01,0x42"
$ws2.Cells.Item(95, 5).Value = "0x19"
$ws2.Cells.Item(95, 6).Value = 400
# Stored as literal text (not a number) in the source data - force text via
# the leading apostrophe, then strip the auto-applied "Text" number format
# so the cell keeps the workbook's default (no explicit style) formatting.
$ws2.Cells.Item(95, 7).Value = "'568631262647113771663628"
$ws2.Cells.Item(95, 7).ClearFormats()
$ws2.Cells.Item(95, 8).Value = 322
$ws2.Cells.Item(95, 9).Value = 25

# ---------------------------------------------------------------
# Sheet 3 "ROW11-FE-LIFTER" - append new row 93
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(93, 1).Value = 45770.82959099537
$ws3.Cells.Item(93, 1).NumberFormat = $dateFmt
$ws3.Cells.Item(93, 2).Value = "0x01,0x90"
$ws3.Cells.Item(93, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item(93, 4).Value = "0x01,0x3e"
$ws3.Cells.Item(93, 5).Value = "0x14"
$ws3.Cells.Item(93, 6).Value = 400
$ws3.Cells.Item(93, 7).Value = [double]"5.68631262647114e+23"
$ws3.Cells.Item(93, 8).Value = 318
$ws3.Cells.Item(93, 9).Value = 20

# ---------------------------------------------------------------
# Sheet 4 "ROW11-MID-LIFTER" - append new row 93
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(93, 1).Value = 45770.95209197917
$ws4.Cells.Item(93, 1).NumberFormat = $dateFmt
$ws4.Cells.Item(93, 2).Value = "0x01,0x90"
$ws4.Cells.Item(93, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item(93, 4).Value = "0x01,0x42"
$ws4.Cells.Item(93, 5).Value = "0x19"
$ws4.Cells.Item(93, 6).Value = 400
$ws4.Cells.Item(93, 7).Value = [double]"5.68631262647114e+23"
$ws4.Cells.Item(93, 8).Value = 322
$ws4.Cells.Item(93, 9).Value = 25
